# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G (K = strikeouts) values for rows 2-17 with the newly calculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 3
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
